$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text so values like "0.9982" or
# "30.096.55" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.096.55"
$ws.Range("E2").Value = "  +6.73%  "
$ws.Range("D3").Value = "1.891.87"
$ws.Range("E3").Value = "  +5.87%  "
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "249.41"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "0.9981"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.5012"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.2872"
$ws.Range("E8").Value = "  +6.87%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06571"
$ws.Range("E9").Value = "  +4.47%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.880.68"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "17.30"
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07245"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.6708"
$ws.Range("E13").Value = "  +6.73%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "85.06"
$ws.Range("E14").Value = "  +6.29%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.836"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.094.98"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "0.9985"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "12.94"
$ws.Range("E18").Value = "  +7.36%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007550"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("B20").Value = "BinanceUSD"
$ws.Range("C20").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D20").Value = "0.9979"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.122.55"
$ws.Range("E21").Value = "  +5.59%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.785"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "5.571"
$ws.Range("E23").Value = "  +6.16%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "9.050"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("B26").Value = "BitcoinCash"
$ws.Range("C26").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D26").Value = "134.92"
$ws.Range("E26").Value = "  +22.98%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "16.80"
$ws.Range("E27").Value = "  +6.52%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "1.955"
$ws.Range("E28").Value = "  +5.36%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.370"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "4.213"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08673"
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.955"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05043"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.145"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.6969"
$ws.Range("E35").Value = "  +6.54%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.691"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "2.296"
$ws.Range("E37").Value = "  +11.76%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.763"
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.9642"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01636"
$ws.Range("E40").Value = "  +5.81%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.058"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "104.71"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.9979"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4229"
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "7.490"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1259"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05660"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "32.54"
$ws.Range("E48").Value = "  +5.82%  "
$ws.Range("D49").Value = "8.277"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "0.3724"
$ws.Range("E50").Value = "  +6.64%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.340"
$ws.Range("E51").Value = "  +3.41%  "
